# Applies the "3e version avec organisation fichiers" update to
# StructureDefinition-Specialite.xlsx:
#   - Metadata sheet: bump Date, change Base Definition URL (Specialite now
#     derives from SavoirFaire instead of the raw FHIR Base type).
#   - Elements sheet: the element that used to be "Specialite.specialite" at
#     row 3 becomes "Specialite.typeSavoirFaire" (inherited from SavoirFaire),
#     and three new rows are introduced below it: Specialite.dateReconnaissance,
#     Specialite.dateAbandon, and a (re-added) Specialite.specialite.

$wb   = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$els  = $wb.Worksheets.Item("Elements")

# ---------------------------------------------------------------------------
# Metadata sheet
# ---------------------------------------------------------------------------
$meta.Range("B8").Value  = "2025-07-21T11:52:46+00:00"
$meta.Range("B18").Value = "https://interop.esante.gouv.fr/ig/fhir/mos/StructureDefinition/SavoirFaire"

# ---------------------------------------------------------------------------
# Elements sheet
# ---------------------------------------------------------------------------

# Row 3 used to describe "Specialite.specialite"; it now describes the
# (renamed/relocated) "Specialite.typeSavoirFaire" element instead.
$els.Range("A3").Value  = "Specialite.typeSavoirFaire"
$els.Range("B3").Value  = "Specialite.typeSavoirFaire"
$shortDef3 = " Le type de savoir-faire (qualifications/autres attributions) d" + [char]0x00E9 + "signe par exemple:** une sp" + [char]0x00E9 + "cialit" + [char]0x00E9 + " ordinale (S);** une comp" + [char]0x00E9 + "tence (C);** etc."
$els.Range("L3").Value  = $shortDef3
$els.Range("M3").Value  = $shortDef3
$els.Range("Z3").Value  = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R04-TypeSavoirFaire?vs"
$els.Range("AF3").Value = "SavoirFaire.typeSavoirFaire"

# Insert three new rows (4,5,6) below it, seeded from the row-2/row-3 layouts
# (border + wrap-text style, and the blank-vs-empty-text cell pattern) so the
# new rows look like the existing ones.
$els.Rows.Item(4).Insert()
$els.Rows.Item(4).Insert()
$els.Rows.Item(4).Insert()

$els.Range("A2:AJ2").Copy($els.Range("A4:AJ4"))
$els.Range("A2:AJ2").Copy($els.Range("A5:AJ5"))
$els.Range("A3:AJ3").Copy($els.Range("A6:AJ6"))

# Columns that must hold an *empty text* value (as opposed to a truly blank
# cell) on a "leaf element" row (the pattern used by row 2 / rows 4 & 5) and
# on a "bound element" row (the pattern used by row 3 / row 6).
$leafEmptyTextCols  = @("D","H","I","J","P","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AI","AJ")
$boundEmptyTextCols = @("D","H","I","J","P","R","S","T","U","V","W","AA","AB","AC","AD","AE","AI","AJ")

foreach ($col in $leafEmptyTextCols) {
    $els.Range($col + "4").Value = "'"
    $els.Range($col + "5").Value = "'"
}
foreach ($col in $boundEmptyTextCols) {
    $els.Range($col + "6").Value = "'"
}

# --- Row 4: Specialite.dateReconnaissance ---
$els.Range("A4").Value  = "Specialite.dateReconnaissance"
$els.Range("B4").Value  = "Specialite.dateReconnaissance"
$els.Range("F4").Value  = "'0"
$els.Range("G4").Value  = "'1"
$els.Range("K4").Value  = "date`n"
$reconnaissanceDef = " Date " + [char]0x00E0 + " laquelle, l" + [char]0x2019 + "organisme donnant l" + [char]0x2019 + "autorisation d" + [char]0x2019 + "exercer une qualification a reconnu cette qualification ou date " + [char]0x00E0 + " laquelle l'attribution a " + [char]0x00E9 + "t" + [char]0x00E9 + " donn" + [char]0x00E9 + "e au professionnel."
$els.Range("L4").Value  = $reconnaissanceDef
$els.Range("M4").Value  = $reconnaissanceDef
$els.Range("AF4").Value = "SavoirFaire.dateReconnaissance"
$els.Range("AG4").Value = "'0"
$els.Range("AH4").Value = "'1"

# --- Row 5: Specialite.dateAbandon ---
$els.Range("A5").Value  = "Specialite.dateAbandon"
$els.Range("B5").Value  = "Specialite.dateAbandon"
$els.Range("F5").Value  = "'0"
$els.Range("G5").Value  = "'1"
$els.Range("K5").Value  = "date`n"
$abandonDef = " Date " + [char]0x00E0 + " laquelle le professionnel a d" + [char]0x00E9 + "clar" + [char]0x00E9 + " renoncer " + [char]0x00E0 + " l" + [char]0x2019 + "exercice d" + [char]0x2019 + "un savoir-faire ou date " + [char]0x00E0 + " laquelle il ne souhaite plus le faire appara" + [char]0x00EE + "tre."
$els.Range("L5").Value  = $abandonDef
$els.Range("M5").Value  = $abandonDef
$els.Range("AF5").Value = "SavoirFaire.dateAbandon"
$els.Range("AG5").Value = "'0"
$els.Range("AH5").Value = "'1"

# --- Row 6: Specialite.specialite (the original element, re-added) ---
$els.Range("A6").Value  = "Specialite.specialite"
$els.Range("B6").Value  = "Specialite.specialite"
$els.Range("F6").Value  = "'0"
$els.Range("G6").Value  = "'1"
$els.Range("K6").Value  = "Coding`n"
$specialiteDef = " Sp" + [char]0x00E9 + "cialit" + [char]0x00E9 + " ordinale."
$els.Range("L6").Value  = $specialiteDef
$els.Range("M6").Value  = $specialiteDef
$els.Range("X6").Value  = "preferred"
$els.Range("Z6").Value  = "https://interop.esante.gouv.fr/ig/fhir/mos/ValueSet/specialite-vs"
$els.Range("AF6").Value = "Specialite.specialite"
$els.Range("AG6").Value = "'0"
$els.Range("AH6").Value = "'1"
